# Update NATMI LR-pair output sheet with newly computed TPM-based values
# and remove the two rows that corresponded to the "ECs" sending cluster
# (rows 6 and 7 in the old layout), since the new TPM pipeline no longer
# emits those sending-cluster rows in this particular pairwise file.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the last two data rows (old rows 6 and 7, the "ECs" sender rows).
$ws.Rows.Item(6).Resize(2).Delete() | Out-Null

# ---- Row 2: FAPs -> Efnb3/Ephb1 -> ECs ----
$ws.Range("A2").Value = "FAPs"
$ws.Range("B2").Value = "Efnb3"
$ws.Range("C2").Value = "Ephb1"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 2
$ws.Range("F2").Value = 0.6666666666666666
$ws.Range("G2").Value = 0.1498043333333333
$ws.Range("H2").Value = 0.449413
$ws.Range("I2").Value = 0.08722868471333377
$ws.Range("J2").Value = 0.08722868471333377
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 1.475208
$ws.Range("N2").Value = 4.425624
$ws.Range("O2").Value = 0.8210007041987012
$ws.Range("P2").Value = 0.8210007041987013
$ws.Range("Q2").Value = 0.220992550968
$ws.Range("R2").Value = 1.988932958712
$ws.Range("S2").Value = 0.07161481157597351
$ws.Range("T2").Value = 0.07161481157597352

# ---- Row 3: FAPs -> Efnb3/Ephb1 -> MuSCs ----
$ws.Range("A3").Value = "FAPs"
$ws.Range("B3").Value = "Efnb3"
$ws.Range("C3").Value = "Ephb1"
$ws.Range("D3").Value = "MuSCs"
$ws.Range("E3").Value = 2
$ws.Range("F3").Value = 0.6666666666666666
$ws.Range("G3").Value = 0.1498043333333333
$ws.Range("H3").Value = 0.449413
$ws.Range("I3").Value = 0.08722868471333377
$ws.Range("J3").Value = 0.08722868471333377
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 0.3216333333333333
$ws.Range("N3").Value = 0.9649
$ws.Range("O3").Value = 0.1789992958012987
$ws.Range("P3").Value = 0.1789992958012987
$ws.Range("Q3").Value = 0.04818206707777778
$ws.Range("R3").Value = 0.4336386037
$ws.Range("S3").Value = 0.01561387313736025
$ws.Range("T3").Value = 0.01561387313736026

# ---- Row 4: MuSCs -> Efnb3/Ephb1 -> ECs ----
$ws.Range("A4").Value = "MuSCs"
$ws.Range("B4").Value = "Efnb3"
$ws.Range("C4").Value = "Ephb1"
$ws.Range("D4").Value = "ECs"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 1.567570333333333
$ws.Range("H4").Value = 4.702711
$ws.Range("I4").Value = 0.9127713152866662
$ws.Range("J4").Value = 0.9127713152866662
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 1.475208
$ws.Range("N4").Value = 4.425624
$ws.Range("O4").Value = 0.8210007041987012
$ws.Range("P4").Value = 0.8210007041987013
$ws.Range("Q4").Value = 2.312492296296
$ws.Range("R4").Value = 20.812430666664
$ws.Range("S4").Value = 0.7493858926227277
$ws.Range("T4").Value = 0.7493858926227278

# ---- Row 5: MuSCs -> Efnb3/Ephb1 -> MuSCs ----
$ws.Range("A5").Value = "MuSCs"
$ws.Range("B5").Value = "Efnb3"
$ws.Range("C5").Value = "Ephb1"
$ws.Range("D5").Value = "MuSCs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 1.567570333333333
$ws.Range("H5").Value = 4.702711
$ws.Range("I5").Value = 0.9127713152866662
$ws.Range("J5").Value = 0.9127713152866662
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 0.3216333333333333
$ws.Range("N5").Value = 0.9649
$ws.Range("O5").Value = 0.1789992958012987
$ws.Range("P5").Value = 0.1789992958012987
$ws.Range("Q5").Value = 0.5041828715444444
$ws.Range("R5").Value = 4.5376458439
$ws.Range("S5").Value = 0.1633854226639384
$ws.Range("T5").Value = 0.1633854226639385
